# Append: 2025-12-16 18:28 JST
#
# A fresh scrape ran at 2025-12-16 18:28:50 JST. It re-fetched every
# already-known listing (so every row's "fetched at" timestamp in column A
# moves forward) and found one brand-new listing, which is inserted as the
# new row 13 - pushing the two rows that used to be 13/14 down to 14/15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-12-16 18:28:50"

# --- Shift the last two existing listings down by one row -----------------
# (row 14 <- old row 13, row 15 <- old row 14), reading the old values out
# first so overwriting row 13 below doesn't clobber them.
$row13Title    = $ws.Range("B13").Value()
$row13Category = $ws.Range("C13").Value()
$row13Price    = $ws.Range("D13").Value()
$row13Deadline = $ws.Range("E13").Value()
$row13Url      = $ws.Range("F13").Value()
$row13Score    = $ws.Range("G13").Value()

$row14Title    = $ws.Range("B14").Value()
$row14Category = $ws.Range("C14").Value()
$row14Price    = $ws.Range("D14").Value()
$row14Deadline = $ws.Range("E14").Value()
$row14Url      = $ws.Range("F14").Value()
$row14Score    = $ws.Range("G14").Value()

$ws.Range("B15").Value = $row14Title
$ws.Range("C15").Value = $row14Category
$ws.Range("D15").Value = $row14Price
$ws.Range("E15").Value = $row14Deadline
$ws.Range("F15").Value = $row14Url
$ws.Range("G15").Value = $row14Score

$ws.Range("B14").Value = $row13Title
$ws.Range("C14").Value = $row13Category
$ws.Range("D14").Value = $row13Price
$ws.Range("E14").Value = $row13Deadline
$ws.Range("F14").Value = $row13Url
$ws.Range("G14").Value = $row13Score

# --- Write the newly discovered listing into row 13 ------------------------
$ws.Range("B13").Value = "【急募】企業のセキュリティ対策を担うエンジニア募集"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5450345"
$ws.Range("G13").Value = 25

# --- Refresh the "fetched at" timestamp for every row (now 2..15) ----------
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- Rebuild the hyperlinks on column F so every row keeps a working link --
# (Range.Hyperlinks.Delete() clears the whole sheet's collection in this
# runtime, so the simplest correct approach is: clear once, then re-add the
# anchor for every F-column cell that should carry a link. Re-applying the
# "Hyperlink" named style afterwards keeps the cell format identical to the
# pre-existing hyperlink cells instead of growing a duplicate style entry.)
$ws.Range("F2").Hyperlinks.Delete()

$linkCells = @(
    @{ Row = 2;  Url = "https://www.lancers.jp/work/detail/5455098" },
    @{ Row = 3;  Url = "https://www.lancers.jp/work/detail/5454985" },
    @{ Row = 4;  Url = "https://www.lancers.jp/work/detail/5217096" },
    @{ Row = 5;  Url = "https://www.lancers.jp/work/detail/5439921" },
    @{ Row = 6;  Url = "https://www.lancers.jp/work/detail/5455160" },
    @{ Row = 7;  Url = "https://www.lancers.jp/work/detail/5455038" },
    @{ Row = 8;  Url = "https://www.lancers.jp/work/detail/5455251" },
    @{ Row = 9;  Url = "https://www.lancers.jp/work/detail/5455029" },
    @{ Row = 10; Url = "https://www.lancers.jp/work/detail/5454857" },
    @{ Row = 11; Url = "https://www.lancers.jp/work/detail/5455015" },
    @{ Row = 12; Url = "https://www.lancers.jp/work/detail/5455067" },
    @{ Row = 13; Url = "https://www.lancers.jp/work/detail/5450345" },
    @{ Row = 14; Url = "https://www.lancers.jp/work/detail/5341051" },
    @{ Row = 15; Url = "https://www.lancers.jp/work/detail/5437544" }
)

foreach ($link in $linkCells) {
    $cell = $ws.Range("F" + $link.Row)
    $ws.Hyperlinks.Add($cell, $link.Url)
    $cell.Style = "Hyperlink"
}
